$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header strings
$ws.Range("D1").Value = "No of Persons"
$ws.Range("J1").Value = "No of App Users"

# Rename Window Orientation values: W -> West, E -> East, N/NE -> North East
# (order chosen so new shared strings are appended as West, East, North East)
$ws.Range("E2").Value = "West"
$ws.Range("E4").Value = "East"
$ws.Range("E3").Value = "North East"
$ws.Range("E5").Value = "West"
$ws.Range("E6").Value = "West"
$ws.Range("E7").Value = "North East"
$ws.Range("E8").Value = "West"
$ws.Range("E9").Value = "East"
$ws.Range("E10").Value = "East"
$ws.Range("E11").Value = "East"
$ws.Range("E12").Value = "West"
$ws.Range("E13").Value = "West"

# Fill empty "No of App Users" cells with 0
$ws.Range("J4").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("J13").Value = 0

# Widen column J to fit new header text (to match bestFit recalculation)
$ws.Columns.Item(10).ColumnWidth = 13.43

# Update selection
$ws.Range("D17").Select() | Out-Null
